$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C, shifting "Link" from C1 to D1.
# The inserted column inherits formatting from the left (B), so the new
# C1 picks up the bordered/bold/centered style automatically.
$ws.Range("C1").EntireColumn.Insert()

# Populate the new column header.
$ws.Range("C1").Value = "Blank"

# D1 (formerly C1, which had no explicit style) needs the same style as
# the other header cells; copy formats over from the now-styled C1.
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)  # xlPasteFormats

# Update the active selection to match the saved view state.
$ws.Range("I9").Select()
